$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember a cell that still carries the original "hyperlink" style (index 2)
# before anything below disturbs it, so it can be restored afterwards.
$ws.Range("M1").Value = "fmt"
$ws.Range("C2").Copy()
$ws.Range("M1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# The data that used to live in row 3 needs to end up in row 4, leaving row 3
# completely empty, and the DNI in column G of that (new) row 4 needs to
# change from the reused value ("90500084Y") to a fresh one ("90500081Y").
# Cut (rather than Insert) moves the row's values/styles down without
# disturbing the worksheet's Hyperlinks collection bookkeeping.
$ws.Range("A3:I3").Cut($ws.Range("A4"))
$ws.Rows.Item(3).ClearFormats()
$ws.Rows.Item(3).ClearContents()

# New DNI value for the relocated row.
$ws.Range("G4").Value = "90500081Y"

# The mailto hyperlink that used to sit on C3 needs to move to C4. The
# Hyperlinks collection doesn't follow cells that get moved/cleared, so
# rebuild it explicitly (C2's has to be re-added too, since clearing the
# collection clears hyperlinks for the whole sheet, not just one range).
$mailAddress = "mailto:juan@example.com"
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), $mailAddress) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), $mailAddress) | Out-Null

# Hyperlinks.Add reformats the target cells with a fresh style; put the
# original hyperlink-cell style (index 2) back on both.
$ws.Range("M1").Copy()
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").Clear()

# Match the view left behind by the edit.
$ws.Range("G4").Select()
